$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $statusText
$ovw.Range("F2").Value = $statusText
$ovw.Range("E3").Value = $statusText
$ovw.Range("F3").Value = $statusText
$ovw.Range("E1").ColumnWidth = 29.166666666666668
$ovw.Range("F1").ColumnWidth = 29.166666666666668

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText
$zh.Range("C1").ColumnWidth = 29.166666666666668
$zh.Range("I1").ColumnWidth = 39.166666666666664
$zh.Range("J1").ColumnWidth = 39.166666666666664

$zh.Range("I2").Value = "573f8fe1-0671-4552-bfd8-037bdf8d7374.md"
$zh.Range("J2").Value = "573f8fe1-0671-4552-bfd8-037bdf8d7374.e7a6bd797cab50c2f65adc9827b7c04bdeeccde0.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-29 21:16:19"

$zh.Range("I3").Value = "ec2c564e-5368-466f-bc81-74a9e2afe9c5.md"
$zh.Range("J3").Value = "ec2c564e-5368-466f-bc81-74a9e2afe9c5.92cde374c60cfa040c8046cb03d2a078fd95248b.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-29 21:16:19"

$zhLink2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d24a5bfd933da5b3ec6a8d7988bb9b32b3b1a6/e2e/573f8fe1-0671-4552-bfd8-037bdf8d7374.md"
$zhLink3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d24a5bfd933da5b3ec6a8d7988bb9b32b3b1a6/e2e/ec2c564e-5368-466f-bc81-74a9e2afe9c5.md"

$zh.Hyperlinks.Delete()
$null = $zh.Hyperlinks.Add($zh.Range("A2"), $zhLink2, "", "", "573f8fe1-0671-4552-bfd8-037bdf8d7374.md")
$null = $zh.Hyperlinks.Add($zh.Range("I2"), $zhLink2, "", "", "573f8fe1-0671-4552-bfd8-037bdf8d7374.md")
$null = $zh.Hyperlinks.Add($zh.Range("A3"), $zhLink3, "", "", "ec2c564e-5368-466f-bc81-74a9e2afe9c5.md")
$null = $zh.Hyperlinks.Add($zh.Range("I3"), $zhLink3, "", "", "ec2c564e-5368-466f-bc81-74a9e2afe9c5.md")

$zh.Range("A2").Style = "HyperLink"
$zh.Range("A3").Style = "HyperLink"
$zh.Range("I2").Style = "HyperLink"
$zh.Range("I3").Style = "HyperLink"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText
$de.Range("C1").ColumnWidth = 29.166666666666668
$de.Range("I1").ColumnWidth = 39.166666666666664
$de.Range("J1").ColumnWidth = 39.166666666666664

$de.Range("I2").Value = "573f8fe1-0671-4552-bfd8-037bdf8d7374.md"
$de.Range("J2").Value = "573f8fe1-0671-4552-bfd8-037bdf8d7374.e7a6bd797cab50c2f65adc9827b7c04bdeeccde0.de-de.xlf"
$de.Range("K2").Value = "2016-08-29 21:16:27"

$de.Range("I3").Value = "ec2c564e-5368-466f-bc81-74a9e2afe9c5.md"
$de.Range("J3").Value = "ec2c564e-5368-466f-bc81-74a9e2afe9c5.92cde374c60cfa040c8046cb03d2a078fd95248b.de-de.xlf"
$de.Range("K3").Value = "2016-08-29 21:16:27"

$deLink2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d24a5bfd933da5b3ec6a8d7988bb9b32b3b1a6/e2e/573f8fe1-0671-4552-bfd8-037bdf8d7374.md"
$deLink3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96d24a5bfd933da5b3ec6a8d7988bb9b32b3b1a6/e2e/ec2c564e-5368-466f-bc81-74a9e2afe9c5.md"

$de.Hyperlinks.Delete()
$null = $de.Hyperlinks.Add($de.Range("A2"), $deLink2, "", "", "573f8fe1-0671-4552-bfd8-037bdf8d7374.md")
$null = $de.Hyperlinks.Add($de.Range("I2"), $deLink2, "", "", "573f8fe1-0671-4552-bfd8-037bdf8d7374.md")
$null = $de.Hyperlinks.Add($de.Range("A3"), $deLink3, "", "", "ec2c564e-5368-466f-bc81-74a9e2afe9c5.md")
$null = $de.Hyperlinks.Add($de.Range("I3"), $deLink3, "", "", "ec2c564e-5368-466f-bc81-74a9e2afe9c5.md")

$de.Range("A2").Style = "HyperLink"
$de.Range("A3").Style = "HyperLink"
$de.Range("I2").Style = "HyperLink"
$de.Range("I3").Style = "HyperLink"
